$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.347.71"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.321.24"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.17"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.92"
$ws.Range("E6").Value = "  -3.07%  "
$ws.Range("E7").Value = "  -2.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.322.00"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.94"
$ws.Range("E11").Value = "  -9.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.137"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.01"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.858.42"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.314.15"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.488.69"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.55"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.963"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "384.35"
$ws.Range("E22").Value = "  +1.70%  "
$ws.Range("E23").Value = "  +5.64%  "
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.20"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.10"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.37"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("E30").Value = "  -2.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "29.04"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "644.87"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.71"
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.52"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.69"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.383"
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0754"
$ws.Range("E40").Value = "  +6.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("E42").Value = "  +13.59%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.996.30"
$ws.Range("E43").Value = "  +3.39%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.127"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  +4.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0402"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.12"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.52"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.126"
$ws.Range("E51").Value = "  -1.25%  "
